$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the stray single-space value in D5 (maintenance type column was left
# with a placeholder space; the fixed template leaves it truly blank).
$ws.Range("D5").Value = ""

# Remove the trailing blank formatting-only row (row 8) that served no
# purpose in the template.
$ws.Rows(8).Delete()
